$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 57612
$ws.Range("B3").Value = 129307
$ws.Range("B4").Value = 119323
$ws.Range("B5").Value = 12712
$ws.Range("B6").Value = 24122
$ws.Range("B7").Value = 1213
$ws.Range("B8").Value = 2658628
$ws.Range("B9").Value = 207103
$ws.Range("B10").Value = 29499
$ws.Range("B11").Value = 589299
$ws.Range("B12").Value = 296374
$ws.Range("B13").Value = 9634
$ws.Range("B14").Value = 162089
$ws.Range("B15").Value = 711779
$ws.Range("B16").Value = 3765
$ws.Range("B17").Value = 341539
$ws.Range("B18").Value = 943213
$ws.Range("B19").Value = 12538
$ws.Range("B21").Value = 941
$ws.Range("B22").Value = 287360
$ws.Range("B23").Value = 188994
$ws.Range("B25").Value = 13832455
$ws.Range("B27").Value = 384887
$ws.Range("B28").Value = 13083
$ws.Range("B29").Value = 142617
$ws.Range("B30").Value = 3477
$ws.Range("B31").Value = 19975
$ws.Range("B32").Value = 5480
$ws.Range("B34").Value = 1113837
$ws.Range("B35").Value = 5728
$ws.Range("B36").Value = 4686
$ws.Range("B37").Value = 1109311
$ws.Range("B38").Value = 102167
$ws.Range("B39").Value = 2619422
$ws.Range("B40").Value = 3815
$ws.Range("B42").Value = 28859
$ws.Range("B43").Value = 228577
$ws.Range("B44").Value = 45474
$ws.Range("B45").Value = 303598
$ws.Range("B46").Value = 91448
$ws.Range("B47").Value = 54621
$ws.Range("B48").Value = 1597103
$ws.Range("B49").Value = 241700
$ws.Range("B50").Value = 10385
$ws.Range("B52").Value = 260133
$ws.Range("B53").Value = 355964
$ws.Range("B54").Value = 214639
$ws.Range("B55").Value = 67249
$ws.Range("B57").Value = 3491
$ws.Range("B58").Value = 116678
$ws.Range("B59").Value = 18414
$ws.Range("B60").Value = 238527
$ws.Range("B62").Value = 83633
$ws.Range("B63").Value = 5285304
$ws.Range("B65").Value = 5720
$ws.Range("B66").Value = 293321
$ws.Range("B67").Value = 3134108
$ws.Range("B69").Value = 311033
$ws.Range("B71").Value = 210667
$ws.Range("B72").Value = 21392
$ws.Range("B73").Value = 3710
$ws.Range("B74").Value = 11642
$ws.Range("B75").Value = 12876
$ws.Range("B76").Value = 199682
$ws.Range("B77").Value = 742198
$ws.Range("B78").Value = 6286
$ws.Range("B79").Value = 14526609
$ws.Range("B80").Value = 1594722
$ws.Range("B81").Value = 2194133
$ws.Range("B82").Value = 964435
$ws.Range("B83").Value = 242819
$ws.Range("B84").Value = 836936
$ws.Range("B85").Value = 3842079
$ws.Range("B86").Value = 43473
$ws.Range("B87").Value = 526307
$ws.Range("B88").Value = 681870
$ws.Range("B89").Value = 335868
$ws.Range("B90").Value = 150260
$ws.Range("B91").Value = 113444
$ws.Range("B92").Value = 100329
$ws.Range("B93").Value = 254472
$ws.Range("B94").Value = 91374
$ws.Range("B95").Value = 54
$ws.Range("B96").Value = 110343
$ws.Range("B97").Value = 506808
$ws.Range("B100").Value = 171131
$ws.Range("B101").Value = 2790
$ws.Range("B102").Value = 231601
$ws.Range("B103").Value = 64746
$ws.Range("B104").Value = 31039
$ws.Range("B105").Value = 33919
$ws.Range("B106").Value = 370528
$ws.Range("B107").Value = 26021
$ws.Range("B108").Value = 12835
$ws.Range("B109").Value = 29860
$ws.Range("B111").Value = 18103
$ws.Range("B113").Value = 2299939
$ws.Range("B114").Value = 244866
$ws.Range("B115").Value = 2391
$ws.Range("B116").Value = 19672
$ws.Range("B117").Value = 95205
$ws.Range("B118").Value = 504847
$ws.Range("B119").Value = 69067
$ws.Range("B120").Value = 46330
$ws.Range("B121").Value = 282890
$ws.Range("B122").Value = 1411474
$ws.Range("B125").Value = 5116
$ws.Range("B126").Value = 164147
$ws.Range("B127").Value = 106727
$ws.Range("B129").Value = 750158
$ws.Range("B130").Value = 360249
$ws.Range("B131").Value = 9343
$ws.Range("B133").Value = 1681063
$ws.Range("B134").Value = 914971
$ws.Range("B135").Value = 2660088
$ws.Range("B136").Value = 829911
$ws.Range("B137").Value = 194930
$ws.Range("B138").Value = 1023565
$ws.Range("B139").Value = 4631336
$ws.Range("B140").Value = 23812
$ws.Range("B142").Value = 4398
$ws.Range("B143").Value = 1819
$ws.Range("B145").Value = 5016
$ws.Range("B146").Value = 2272
$ws.Range("B147").Value = 403106
$ws.Range("B148").Value = 39664
$ws.Range("B149").Value = 657716
$ws.Range("B150").Value = 4834
$ws.Range("B152").Value = 60769
$ws.Range("B153").Value = 374586
$ws.Range("B154").Value = 230826
$ws.Range("B155").Value = 20
$ws.Range("B157").Value = 1564355
$ws.Range("B158").Value = 10432
$ws.Range("B159").Value = 3407283
$ws.Range("B160").Value = 96186
$ws.Range("B161").Value = 33022
$ws.Range("B162").Value = 9496
$ws.Range("B163").Value = 900138
$ws.Range("B164").Value = 632399
$ws.Range("B165").Value = 20856
$ws.Range("B166").Value = 1070
$ws.Range("B169").Value = 39038
$ws.Range("B170").Value = 1193
$ws.Range("B171").Value = 12391
$ws.Range("B173").Value = 281777
$ws.Range("B174").Value = 4150039
$ws.Range("B175").Value = 31575640
$ws.Range("B176").Value = 41310
$ws.Range("B177").Value = 1974056
$ws.Range("B178").Value = 493266
$ws.Range("B179").Value = 4398903
$ws.Range("B180").Value = 159569
$ws.Range("B181").Value = 86022
$ws.Range("B183").Value = 180609
$ws.Range("B184").Value = 2772
$ws.Range("B185").Value = 278135
$ws.Range("B186").Value = 5715
$ws.Range("B187").Value = 90750
$ws.Range("B188").Value = 37534
